# "Cap nhat file Tasks" - update the % Build (column D) and % Testing
# (column H) progress values on Sheet1 for rows 2-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# % Build column (D) - tasks marked complete become 1 (100%)
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 0

# % Testing column (H) - remains 0 (0%) for every task
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
